$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Title" column for the 5 previously-blank rows
# (row 46 -> C46 .. row 50 -> C50), matching the source material's topics.
$ws.Range("C46").Value = "Threads group, Multithreading and Thread Call"
$ws.Range("C47").Value = "1-4 Only"
$ws.Range("C48").Value = "1-4 Only"
$ws.Range("C49").Value = "Saturday Holiday"
$ws.Range("C50").Value = "Stream Classes, Character Stream, BufferedReader"

# Column C has bestFit/autosize behavior - widen it so the new, longer
# entries are not truncated (mirrors Excel's automatic "best fit" resize
# that happens when longer text is typed into such a column).
$ws.Columns.Item(3).ColumnWidth = 45.3
